# Update crypto price/volume data per the Feb 9 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = '''46.146.60'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +3.81%  '
$c = $ws.Range("D3")
$c.Value = '''2.448.58'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.18%  '
$c = $ws.Range("D4")
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.Value = '''322.13'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.73%  '
$c = $ws.Range("D6")
$c.Value = '''104.69'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.96%  '
$ws.Range("E7").Value = '  +1.25%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +4.99%  '
$ws.Range("E10").Value = '  +2.25%  '
$c = $ws.Range("D11")
$c.Value = '''0.0806'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.98%  '
$ws.Range("E12").Value = '  -2.09%  '
$c = $ws.Range("D13")
$c.Value = '''18.32'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -4.14%  '
$ws.Range("E14").Value = '  +2.08%  '
$c = $ws.Range("D15")
$c.Value = '''2.833.65'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.20%  '
$c = $ws.Range("D16")
$c.Value = '''2.440.06'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("E17").Value = '  +1.28%  '
$c = $ws.Range("D18")
$c.Value = '''46.040.01'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +3.88%  '
$c = $ws.Range("D19")
$c.Value = '''12.66'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.03%  '
$ws.Range("E20").Value = '  +0.89%  '
$c = $ws.Range("D21")
$c.Value = '''0.0₃0932'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.63%  '
$c = $ws.Range("D22")
$c.Value = '''70.96'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +3.42%  '
$ws.Range("E23").Value = '  +4.93%  '
$c = $ws.Range("D24")
$c.Value = '''246.82'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +2.01%  '
$c = $ws.Range("D25")
$c.Value = '''2.51'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.13%  '
# Row 26: coin identity swap plus updated price/volume
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D26")
$c.Value = '''0.999'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.14%  '
# Row 27: coin identity swap plus updated price/volume
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D27")
$c.Value = '''25.87'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +3.05%  '
$c = $ws.Range("D28")
$c.Value = '''2.29'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.72%  '
$c = $ws.Range("D29")
$c.Value = '''9.69'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.52%  '
$c = $ws.Range("D30")
$c.Value = '''34.24'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +3.02%  '
$c = $ws.Range("D31")
$c.Value = '''49.37'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.82%  '
$ws.Range("E32").Value = '  +2.85%  '
$c = $ws.Range("D33")
$c.Value = '''19.82'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.58%  '
$c = $ws.Range("D34")
$c.Value = '''5.35'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +3.68%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E36").Value = '  -0.93%  '
$c = $ws.Range("D37")
$c.Value = '''4.53'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.63%  '
$ws.Range("E38").Value = '  +0.27%  '
$c = $ws.Range("D39")
$c.Value = '''2.95'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.88%  '
$c = $ws.Range("D40")
$c.Value = '''127.94'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +4.60%  '
$ws.Range("E41").Value = '  +2.02%  '
$ws.Range("E42").Value = '  +0.68%  '
$c = $ws.Range("D43")
$c.Value = '''20.80'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("E44").Value = '  +1.56%  '
$c = $ws.Range("D45")
$c.Value = '''1.969.66'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("E47").Value = '  -4.86%  '
$c = $ws.Range("D48")
$c.Value = '''1.86'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +12.99%  '
$ws.Range("E49").Value = '  -3.95%  '
$ws.Range("E50").Value = '  +8.31%  '
$c = $ws.Range("D51")
$c.Value = '''77.86'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +5.15%  '
